# Generate Report for Handback
# The underlying report data was regenerated: file "16f3d321-...md" became
# "95f8a318-...md" and file "3cf24830-...md" became "ffffd0bf626d-...md", with
# new xliff hashes / handback timestamps. Update cell values and hyperlink
# display text on all three sheets to reflect the refreshed report.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "95f8a318-422c-44d2-9ac2-bbf9c37f8ca0.md"
$ws.Range("G2").Value = "2016-09-05 11:32:40"

$ws.Range("A3").Value = "ffffd0bf626d-87c6-4a6e-a325-b983ce8ec19b.md"
$ws.Range("G3").Value = "2016-09-05 11:32:40"

# Rebuild hyperlinks with the refreshed display text (targets are unchanged).
$ov_rid2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/101882d2ecce48aad1a0e290a8b9b1b9a5227c4c/e2e/16f3d321-3073-4df4-ab5a-be6fa124d0fb.md"
$ov_rid3 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/101882d2ecce48aad1a0e290a8b9b1b9a5227c4c/e2e/3cf24830-eb2a-49b9-9013-9fe3ca001af4.md"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), $ov_rid2, "", "", "e2e\95f8a318-422c-44d2-9ac2-bbf9c37f8ca0.md")
$ws.Hyperlinks.Add($ws.Range("B3"), $ov_rid3, "", "", "e2e\ffffd0bf626d-87c6-4a6e-a325-b983ce8ec19b.md")

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "95f8a318-422c-44d2-9ac2-bbf9c37f8ca0.md"
$ws.Range("G2").Value = "95f8a318-422c-44d2-9ac2-bbf9c37f8ca0.65d14dbc9e20dcccae608ac60dac112b371f00a9.zh-cn.xlf"
$ws.Range("H2").Value = "2016-09-05 11:32:35"
$ws.Range("I2").Value = "95f8a318-422c-44d2-9ac2-bbf9c37f8ca0.md"
$ws.Range("J2").Value = "95f8a318-422c-44d2-9ac2-bbf9c37f8ca0.65d14dbc9e20dcccae608ac60dac112b371f00a9.zh-cn.xlf"
$ws.Range("K2").Value = "2016-09-05 11:33:05"

$ws.Range("A3").Value = "ffffd0bf626d-87c6-4a6e-a325-b983ce8ec19b.md"
$ws.Range("G3").Value = "95f8a318-422c-44d2-9ac2-bbf9c37f8ca0.65d14dbc9e20dcccae608ac60dac112b371f00a9.zh-cn.xlf"
$ws.Range("H3").Value = "2016-09-05 11:32:35"
$ws.Range("I3").Value = "ffffd0bf626d-87c6-4a6e-a325-b983ce8ec19b.md"
$ws.Range("J3").Value = "95f8a318-422c-44d2-9ac2-bbf9c37f8ca0.65d14dbc9e20dcccae608ac60dac112b371f00a9.zh-cn.xlf"
$ws.Range("K3").Value = "2016-09-05 11:33:05"

$zh_rid2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/101882d2ecce48aad1a0e290a8b9b1b9a5227c4c/e2e/16f3d321-3073-4df4-ab5a-be6fa124d0fb.md"
$zh_rid3 = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/4bac0136f1a40922f34072ae4872afcbfebbb00c/e2e/16f3d321-3073-4df4-ab5a-be6fa124d0fb.md"
$zh_rid4 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/101882d2ecce48aad1a0e290a8b9b1b9a5227c4c/e2e/3cf24830-eb2a-49b9-9013-9fe3ca001af4.md"
$zh_rid5 = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/4bac0136f1a40922f34072ae4872afcbfebbb00c/e2e/3cf24830-eb2a-49b9-9013-9fe3ca001af4.md"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $zh_rid2, "", "", "95f8a318-422c-44d2-9ac2-bbf9c37f8ca0.md")
$ws.Hyperlinks.Add($ws.Range("I2"), $zh_rid3, "", "", "95f8a318-422c-44d2-9ac2-bbf9c37f8ca0.md")
$ws.Hyperlinks.Add($ws.Range("A3"), $zh_rid4, "", "", "ffffd0bf626d-87c6-4a6e-a325-b983ce8ec19b.md")
$ws.Hyperlinks.Add($ws.Range("I3"), $zh_rid5, "", "", "ffffd0bf626d-87c6-4a6e-a325-b983ce8ec19b.md")

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "95f8a318-422c-44d2-9ac2-bbf9c37f8ca0.md"
$ws.Range("G2").Value = "95f8a318-422c-44d2-9ac2-bbf9c37f8ca0.65d14dbc9e20dcccae608ac60dac112b371f00a9.de-de.xlf"
$ws.Range("H2").Value = "2016-09-05 11:32:40"
$ws.Range("I2").Value = "95f8a318-422c-44d2-9ac2-bbf9c37f8ca0.md"
$ws.Range("J2").Value = "95f8a318-422c-44d2-9ac2-bbf9c37f8ca0.65d14dbc9e20dcccae608ac60dac112b371f00a9.de-de.xlf"
$ws.Range("K2").Value = "2016-09-05 11:33:19"

$ws.Range("A3").Value = "ffffd0bf626d-87c6-4a6e-a325-b983ce8ec19b.md"
$ws.Range("G3").Value = "95f8a318-422c-44d2-9ac2-bbf9c37f8ca0.65d14dbc9e20dcccae608ac60dac112b371f00a9.de-de.xlf"
$ws.Range("H3").Value = "2016-09-05 11:32:40"
$ws.Range("I3").Value = "ffffd0bf626d-87c6-4a6e-a325-b983ce8ec19b.md"
$ws.Range("J3").Value = "95f8a318-422c-44d2-9ac2-bbf9c37f8ca0.65d14dbc9e20dcccae608ac60dac112b371f00a9.de-de.xlf"
$ws.Range("K3").Value = "2016-09-05 11:33:19"

$de_rid2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/101882d2ecce48aad1a0e290a8b9b1b9a5227c4c/e2e/16f3d321-3073-4df4-ab5a-be6fa124d0fb.md"
$de_rid3 = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/cc45595fe8d62f073b76d7cea99fb80dd91673bd/e2e/16f3d321-3073-4df4-ab5a-be6fa124d0fb.md"
$de_rid4 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/101882d2ecce48aad1a0e290a8b9b1b9a5227c4c/e2e/3cf24830-eb2a-49b9-9013-9fe3ca001af4.md"
$de_rid5 = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/cc45595fe8d62f073b76d7cea99fb80dd91673bd/e2e/3cf24830-eb2a-49b9-9013-9fe3ca001af4.md"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $de_rid2, "", "", "95f8a318-422c-44d2-9ac2-bbf9c37f8ca0.md")
$ws.Hyperlinks.Add($ws.Range("I2"), $de_rid3, "", "", "95f8a318-422c-44d2-9ac2-bbf9c37f8ca0.md")
$ws.Hyperlinks.Add($ws.Range("A3"), $de_rid4, "", "", "ffffd0bf626d-87c6-4a6e-a325-b983ce8ec19b.md")
$ws.Hyperlinks.Add($ws.Range("I3"), $de_rid5, "", "", "ffffd0bf626d-87c6-4a6e-a325-b983ce8ec19b.md")
